$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.30%"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "'12"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'41.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'4.78%"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "'12"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'5.670"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.83%"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").Value = "'12"
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").Value = "'0.08299"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.09%"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").Value = "'12"
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = "'2.043"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.63%"
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").Value = "'12"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'8.768"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.45%"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "'12"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'4.531"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.76%"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'12"
$ws.Range("G8").Style = "Normal"
$ws.Range("E9").Value = "'0.55%"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Value = "'12"
$ws.Range("G9").Style = "Normal"
$ws.Range("D10").Value = "'0.9256"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.23%"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'12"
$ws.Range("G10").Style = "Normal"
$ws.Range("D11").Value = "'0.1267"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.50%"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").Value = "'12"
$ws.Range("G11").Style = "Normal"
$ws.Range("D12").Value = "'0.1951"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.62%"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").Value = "'12"
$ws.Range("G12").Style = "Normal"
$ws.Range("D13").Value = "'0.09416"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.46%"
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").Value = "'12"
$ws.Range("G13").Style = "Normal"
$ws.Range("D14").Value = "'0.03886"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'9.04%"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").Value = "'12"
$ws.Range("G14").Style = "Normal"
$ws.Range("E15").Value = "'0.89%"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").Value = "'12"
$ws.Range("G15").Style = "Normal"
$ws.Range("D16").Value = "'0.001310"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.48%"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").Value = "'12"
$ws.Range("G16").Style = "Normal"
$ws.Range("D17").Value = "'0.006102"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-3.62%"
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").Value = "'12"
$ws.Range("G17").Style = "Normal"
$ws.Range("G18").Value = "'12"
$ws.Range("G18").Style = "Normal"
$ws.Range("D19").Value = "'3.439"
$ws.Range("D19").Style = "Normal"
$ws.Range("G19").Value = "'12"
$ws.Range("G19").Style = "Normal"
$ws.Range("E20").Value = "'1.50%"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").Value = "'12"
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").Value = "'8.366"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-4.43%"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").Value = "'12"
$ws.Range("G21").Style = "Normal"
$ws.Range("E22").Value = "'1.94%"
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").Value = "'12"
$ws.Range("G22").Style = "Normal"
$ws.Range("E23").Value = "'-1.22%"
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").Value = "'12"
$ws.Range("G23").Style = "Normal"
$ws.Range("D24").Value = "'0.04407"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.02%"
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").Value = "'12"
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").Value = "'0.001258"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.07%"
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").Value = "'12"
$ws.Range("G25").Style = "Normal"
$ws.Range("E26").Value = "'-6.55%"
$ws.Range("E26").Style = "Normal"
$ws.Range("G26").Value = "'12"
$ws.Range("G26").Style = "Normal"
$ws.Range("E27").Value = "'0.97%"
$ws.Range("E27").Style = "Normal"
$ws.Range("G27").Value = "'12"
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").Value = "'12"
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = "'12"
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = "'12"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "'12"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "'12"
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = "'12"
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = "'12"
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = "'12"
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = "'12"
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = "'12"
$ws.Range("G37").Style = "Normal"
$ws.Range("G38").Value = "'12"
$ws.Range("G38").Style = "Normal"
$ws.Range("D39").Value = "'0.02782"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'11.67%"
$ws.Range("E39").Style = "Normal"
$ws.Range("G39").Value = "'12"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.05508"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.02%"
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").Value = "'12"
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = "'0.007910"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'5.67%"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").Value = "'12"
$ws.Range("G41").Style = "Normal"
$ws.Range("E42").Value = "'0.99%"
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").Value = "'12"
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").Value = "'0.008932"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-9.87%"
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").Value = "'12"
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = "'0.002141"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'1.34%"
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").Value = "'12"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.01190"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'18.06%"
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").Value = "'12"
$ws.Range("G45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006745"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'1.32%"
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").Value = "'12"
$ws.Range("G46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.26%"
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").Value = "'12"
$ws.Range("G47").Style = "Normal"
$ws.Range("D48").Value = "'0.003190"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'5.13%"
$ws.Range("E48").Style = "Normal"
$ws.Range("G48").Value = "'12"
$ws.Range("G48").Style = "Normal"
$ws.Range("E49").Value = "'0.26%"
$ws.Range("E49").Style = "Normal"
$ws.Range("G49").Value = "'12"
$ws.Range("G49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.26%"
$ws.Range("E50").Style = "Normal"
$ws.Range("G50").Value = "'12"
$ws.Range("G50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.26%"
$ws.Range("E51").Style = "Normal"
$ws.Range("G51").Value = "'12"
$ws.Range("G51").Style = "Normal"
